# Update "paises.xlsx" country data: several countries' case counts grew and
# overtook their neighbours in the (descending, by Casos totales) ranking, so
# those rows need to be re-sorted; a handful of other rows get pure value
# refreshes without changing rank. Finally the "last updated" timestamp in A1
# is bumped.
#
# Rather than re-writing the whole sheet, we only touch the rows whose
# contents actually change between the old and new snapshot - for the
# re-sorted blocks that means writing the newly-promoted country's fresh
# numbers into the earlier row and shifting the bumped countries' (unchanged)
# numbers down into the following rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=4;   Name="Estados Unidos";               Vals=@(1676944, 10116, 451176, 1126733, 0, 352, 99035) },
    @{ Row=11;  Name="Alemania";                      Vals=@(180219, 233, 160300, 11548, 0, 5, 8371) },
    @{ Row=41;  Name="Egipto";                        Vals=@(17265, 752, 4807, 11694, 0, 29, 764) },
    @{ Row=42;  Name="Israel";                        Vals=@(16717, 5, 14153, 2285, 0, 0, 279) },
    @{ Row=43;  Name="Japon";                         Vals=@(16536, 0, 13244, 2484, 0, 0, 808) },
    @{ Row=101; Name="Maldivas";                      Vals=@(1371, 58, 144, 1223, 0, 0, 4) },
    @{ Row=103; Name="Sri Lanka";                     Vals=@(1141, 52, 674, 458, 0, 0, 9) },
    @{ Row=115; Name="Costa Rica";                    Vals=@(930, 12, 620, 300, 0, 0, 10) },
    @{ Row=116; Name="Zambia";                        Vals=@(920, 0, 336, 577, 0, 0, 7) },
    @{ Row=152; Name="Mauritania";                    Vals=@(227, 0, 15, 206, 0, 0, 6) },
    @{ Row=173; Name="Comoras";                       Vals=@(87, 9, 21, 65, 0, 0, 1) },
    @{ Row=174; Name="Siria";                         Vals=@(86, 16, 41, 41, 0, 0, 4) },
    @{ Row=175; Name="Malaui";                        Vals=@(82, 0, 28, 50, 0, 1, 4) },
    @{ Row=176; Name="Liechtenstein";                 Vals=@(82, 0, 55, 26, 0, 0, 1) },
    @{ Row=198; Name="Santa Lucia";                   Vals=@(18, 0, 18, 0, 0, 0, 0) },
    @{ Row=199; Name="Nueva Caledonia";                Vals=@(18, 0, 18, 0, 0, 0, 0) },
    @{ Row=200; Name="Belice";                        Vals=@(18, 0, 16, 0, 0, 0, 2) },
    @{ Row=209; Name="Groenlandia";                   Vals=@(11, 0, 11, 0, 0, 0, 0) },
    @{ Row=210; Name="Seychelles";                    Vals=@(11, 0, 11, 0, 0, 0, 0) },
    @{ Row=214; Name="Sahara Occidental";              Vals=@(6, 0, 6, 0, 0, 0, 0) },
    @{ Row=215; Name="Bonaire, San Eustaquio y Saba";  Vals=@(6, 0, 6, 0, 0, 0, 0) },
    @{ Row=216; Name="San Bartolome";                  Vals=@(6, 0, 6, 0, 0, 0, 0) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Name
    $col = 2
    foreach ($v in $u.Vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}

$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 21:05"
